$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.687.23"
$ws.Range("E2").Value = "  -5.95%  "
$ws.Range("D3").Value = "3.295.35"
$ws.Range("E3").Value = "  -6.18%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "556.89"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").Value = "182.82"
$ws.Range("E6").Value = "  -5.30%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("D9").Value = "3.287.12"
$ws.Range("E9").Value = "  -6.08%  "
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  -10.14%  "
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").Value = "  -6.27%  "
$ws.Range("D12").Value = "47.31"
$ws.Range("E12").Value = "  -8.18%  "
$ws.Range("E13").Value = "  -7.19%  "
$ws.Range("D14").Value = "647.75"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "8.63"
$ws.Range("E15").Value = "  -6.05%  "
$ws.Range("D16").Value = "3.822.67"
$ws.Range("E16").Value = "  -6.16%  "
$ws.Range("D17").Value = "18.11"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "65.753.51"
$ws.Range("E18").Value = "  -5.86%  "
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").Value = "3.293.08"
$ws.Range("E20").Value = "  -6.37%  "
$ws.Range("E21").Value = "  -8.52%  "
$ws.Range("D22").Value = "0.904"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "18.11"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "107.41"
$ws.Range("E24").Value = "  +7.80%  "
$ws.Range("D25").Value = "4.88"
$ws.Range("E25").Value = "  -9.00%  "
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").Value = "  -7.61%  "
$ws.Range("E27").Value = "  -7.57%  "
$ws.Range("D28").Value = "9.58"
$ws.Range("E28").Value = "  -5.52%  "
$ws.Range("E29").Value = "  -7.82%  "
$ws.Range("D30").Value = "30.24"
$ws.Range("E30").Value = "  -7.83%  "
$ws.Range("D31").Value = "3.89"
$ws.Range("E31").Value = "  -8.63%  "
$ws.Range("D32").Value = "6.27"
$ws.Range("E32").Value = "  -7.46%  "
$ws.Range("D33").Value = "11.03"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("E34").Value = "  -5.21%  "
$ws.Range("D35").Value = "3.787.02"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "57.36"
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "523.09"
$ws.Range("E38").Value = "  -9.44%  "
$ws.Range("D39").Value = "0.0₃0733"
$ws.Range("E39").Value = "  -8.00%  "
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  -7.63%  "
$ws.Range("D41").Value = "0.129"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  -6.61%  "
$ws.Range("D43").Value = "32.89"
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  -10.45%  "
$ws.Range("E45").Value = "  -10.30%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.26"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0414"
$ws.Range("E47").Value = "  -7.08%  "
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("E49").Value = "  -9.74%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "1.27"
$ws.Range("E51").Value = "  +1.80%  "
